# Update the indicator metadata values on the "Пример" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10: Organization website
$ws.Range("B10").Value = "www.stat.gov.kg"

# B9: Contact phone number
$ws.Range("B9").Value = "0(312) 32 55 46"

# B4: Indicator name/description text
$ws.Range("B4").Value = "3.7.2 Показатель рождаемости среди девушек-подростков ( в возрасте от 15 до 19 лет) на 1000 девушек-подростков в той же возрастной группе"

# Move the active selection to B4 to match the saved view state
$ws.Range("B4").Select()
